# 自动更新价格数据：在表头下方插入最新一天的价格记录，
# 其余历史数据行整体下移一行。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 在第 2 行（表头之后、最早插入的数据行之前）插入一整行，
# 把原来的第 2..47 行数据依次下移到第 3..48 行。
$ws.Range("A2:D2").Insert()

# 新插入的行会继承上方表头的格式（粗体/边框等），清除格式
# 使其与其余普通数据行保持一致（无显式样式）。
$ws.Range("A2:D2").ClearFormats()

# 日期列在原表中以纯文本形式存储（如 "2026-01-05"），
# 先将数字格式设为文本，避免 "2026-01-06" 被自动识别成日期序列值。
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-06"

# 写入文本值后再次清除格式，去掉刚才临时设置的文本数字格式，
# 使该单元格与其它数据行一样不带任何显式样式。
$ws.Range("A2:D2").ClearFormats()

# 新一天的价格数据（与近期数值保持一致）。
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
